$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A3 changes style from s="2" to s="3" ---
# Copy the format of an existing s="3" cell (A4) onto A3, leaving its
# value/hyperlink untouched.
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)

# --- New rows 20 and 21: B2C hotel booking flow URLs ---
$ws.Range("A20").Value = "https://www.fabhotels.com/uiuiuiiuyad"
$ws.Range("A21").Value = "https://www.fabhotels.com/uiuiuiiuyad/hihuiu"

# Hyperlinks for the two new rows (added first so the style copy below
# is the final, authoritative formatting applied to these cells).
$ws.Hyperlinks.Add($ws.Range("A20"), "https://www.fabhotels.com/uiuiuiiuyad")
$ws.Hyperlinks.Add($ws.Range("A21"), "https://www.fabhotels.com/uiuiuiiuyad/hihuiu")

# Adding hyperlinks auto-registers a built-in "Hyperlink" named cell style;
# drop it again since it is not used in the source workbook/style sheet.
$wb.Styles.Item("Hyperlink").Delete()

# Match existing s="2" style (same as A2/old A3) on the new rows.
$ws.Range("A2").Copy()
$ws.Range("A20:A21").PasteSpecial(-4122)

Write-Host "edit applied"
